$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update ISI2 quantities (column F) for rows 2-5
$ws.Range("F2").Value = 50
$ws.Range("F3").Value = 20
$ws.Range("F4").Value = 30
$ws.Range("F5").Value = 10

# Clear the third unit (SATUAN3/HARGAJUAL3/BARCODE3) for row 3, and zero out ISI3
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = ""
$ws.Range("N3").Value = ""
$ws.Range("O3").Value = ""

# Update sheet view: scroll so column B is the left-most visible column,
# and move the active selection to N11
$ws.Application.ActiveWindow.ScrollColumn = 2
$ws.Range("N11").Select()
